$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "complete" / "completed" / "complété" / "RI_593"
$ws.Range("A10").Value = "complete"
$ws.Range("B10").Value = "completed"
$ws.Range("C10").Value = "complété"
$ws.Range("D10").Value = "RI_593"

# Row 11: "In work" / "onGoing" / "enContinue" / "RI_596"
$ws.Range("A11").Value = "In work"
$ws.Range("B11").Value = "onGoing"
$ws.Range("C11").Value = "enContinue"
$ws.Range("D11").Value = "RI_596"

# Match the author's final selection position
$ws.Range("A15").Select() | Out-Null
